$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four task descriptions in column E to reflect that the
# personnel modules now also cover "phòng ban" (department).
$ws.Range("E8").Value = "Xây dựng module sửa nhân sự và phòng ban"
$ws.Range("E9").Value = "Xây dựng module xóa nhân sự và phòng ban"
$ws.Range("E10").Value = "Xây dựng module tìm kiếm nhân sự và phòng ban"
$ws.Range("E7").Value = "Xây dựng module thêm nhân sự và pb mới"

# The longer text in E10 now wraps onto a second line; match the row's
# auto-computed height for that two-line wrap.
$ws.Rows.Item(10).RowHeight = 31.5

# Restore the active selection to E7.
$ws.Range("E7").Select()
